$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data (rows 2..23), columns A (id), B (BFt), C (FuelMix)
$data = @(
    @("id_DK_Central_BH_Biogas", "Biogas", 1.5),
    @("id_DK_Central_BP_Biogas", "Biogas", 2.595057034220532),
    @("id_DK_Central_IndustryH_Biogas", "Biogas", 0.6578947368421052),
    @("id_DK_Central_BH_Biomass", "Biomass", 1.014458272327965),
    @("id_DK_Central_BP_Biomass", "Biomass", 4.261060393258427),
    @("id_DK_Central_IndustryH_Biomass", "Biomass", 0.9711538461538463),
    @("id_DK_Central_BP_Coal", "Coal", 2.650371944739639),
    @("id_DK_Central_BH_Natgas", "Natgas", 1.019108280254777),
    @("id_DK_Central_BP_Natgas", "Natgas", 2.489051094890511),
    @("id_DK_Central_IndustryH_Natgas", "Natgas", 1),
    @("id_DK_Central_BH_Oil", "Oil", 1.127044025157233),
    @("id_DK_Central_BP_Oil", "Oil", 3.369565217391304),
    @("id_DK_Central_IndustryH_Oil", "Oil", 0.5000000000000001),
    @("id_DK_Central_BH_Waste", "Waste", 1.035805626598465),
    @("id_DK_Central_BP_Waste", "Waste", 6.546302050963331),
    @("id_DK_nan_CD_Biogas", "Biogas", 2.846153846153846),
    @("id_DK_nan_IndustryE_Biogas", "Biogas", 2.110204081632653),
    @("id_DK_nan_IndustryE_Biomass", "Biomass", "inf"),
    @("id_DK_nan_CD_Coal", "Coal", 3.210526315789474),
    @("id_DK_nan_IndustryE_Natgas", "Natgas", 2.452631578947368),
    @("id_DK_nan_CD_Oil", "Oil", 2.6),
    @("id_DK_nan_IndustryE_Oil", "Oil", 11)
)

# Delete the rows that no longer exist (old rows 24..35), since the new
# table only spans down to row 23.
$oldLastRow = 35
$newLastRow = 1 + $data.Count
if ($oldLastRow -gt $newLastRow) {
    $ws.Range("A$($newLastRow + 1):C$oldLastRow").EntireRow.Delete() | Out-Null
}

# Write out the rows 2..23 with the final id/BFt/FuelMix values.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
